# Apply the "keep add case for sdl-6761" change:
#  - Update the listGraphNames sheet's "response" column to "responseData"
#    and reformat its sample value from a bracketed/quoted list to a plain
#    comma separated list.
#  - Move listGraphNames so it sits after deleteRelations (end of the
#    relation sheets block) instead of before getRelationById.
#  - Add two new sheets: createInstanceGraph and generateKg.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the existing listGraphNames sheet content.
# ---------------------------------------------------------------------
$listGraphNames = $wb.Worksheets.Item("listGraphNames")
$listGraphNames.Range("C1").Value = "responseData"
$listGraphNames.Range("C2").Value = "janusgraph_iot_demo_dev_kg,janusgraph_iot_demo_dev_instance_kg"

# ---------------------------------------------------------------------
# 2) Build "createInstanceGraph" from a copy of listGraphNames so it
#    inherits the same header/value styling, then insert the extra
#    "graphName" column and fill in the new values.
# ---------------------------------------------------------------------
$listGraphNames.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$createInstanceGraph = $wb.Worksheets.Item($wb.Worksheets.Count)
$createInstanceGraph.Name = "createInstanceGraph"

# Insert the "graphName" column between description (B) and responseData (C).
$createInstanceGraph.Columns.Item(3).Insert()

$createInstanceGraph.Range("A1").Value = "test-id"
$createInstanceGraph.Range("B1").Value = "description"
$createInstanceGraph.Range("C1").Value = "graphName"
$createInstanceGraph.Range("D1").Value = "responseData"
$createInstanceGraph.Range("E1").Value = "rspStatus"
$createInstanceGraph.Range("F1").Value = "rspCode"
$createInstanceGraph.Range("G1").Value = "rspMessage"

$createInstanceGraph.Range("A2").Value = "iot-lpg-create-instance-graph"
$createInstanceGraph.Range("B2").Value = "create instance graph"
$createInstanceGraph.Range("C2").Value = "test6761"
$createInstanceGraph.Range("D2").Value = "janusgraph_iot_demo_dev_kg,test6761,janusgraph_iot_demo_dev_instance_kg"
$createInstanceGraph.Range("E2").Value = 200
$createInstanceGraph.Range("F2").Value = 100000
$createInstanceGraph.Range("G2").Value = "OK"

$createInstanceGraph.Range("C2").WrapText = $true
$createInstanceGraph.Range("D2").WrapText = $true
$createInstanceGraph.Columns.AutoFit()

# ---------------------------------------------------------------------
# 3) Build "generateKg" the same way, with graphName/entityLabels columns
#    inserted after description (the old "response" column slides right
#    and is reused/overwritten to hold the graphql query value).
# ---------------------------------------------------------------------
$listGraphNames.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$generateKg = $wb.Worksheets.Item($wb.Worksheets.Count)
$generateKg.Name = "generateKg"

$generateKg.Range("C1:D1").EntireColumn.Insert()

$generateKg.Range("A1").Value = "test-id"
$generateKg.Range("B1").Value = "description"
$generateKg.Range("C1").Value = "graphName"
$generateKg.Range("D1").Value = "entityLabels"
$generateKg.Range("E1").Value = "graphql"
$generateKg.Range("F1").Value = "rspStatus"
$generateKg.Range("G1").Value = "rspCode"
$generateKg.Range("H1").Value = "rspMessage"

$generateKg.Range("A2").Value = "iot-lpg-generate-kg"
$generateKg.Range("B2").Value = "generate kg"
$generateKg.Range("C2").Value = "test6761"
$generateKg.Range("D2").Value = "sdl_b"
$generateKg.Range("E2").Value = "{`n    sdl_b (cond: ""`", authInfo: ""`", order: ""`") {`n        id`n        hobby`n    }`n}"
$generateKg.Range("F2").Value = 200
$generateKg.Range("G2").Value = 100000
$generateKg.Range("H2").Value = "OK"

$generateKg.Range("C2").WrapText = $true
$generateKg.Range("D2").WrapText = $true
$generateKg.Range("E2").WrapText = $true
$generateKg.Columns.AutoFit()

# ---------------------------------------------------------------------
# 4) Move listGraphNames so it now comes after deleteRelations (right
#    before the two newly added sheets).
# ---------------------------------------------------------------------
$deleteRelations = $wb.Worksheets.Item("deleteRelations")
$listGraphNames.Move($null, $deleteRelations)

# listGraphNames stays the active/selected tab after the reshuffle. Re-fetch
# it by name first since the old COM reference's cached index goes stale
# across a Move().
$wb.Worksheets.Item("listGraphNames").Activate()
